$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("哈希")

# Copy formatting from row 12 down to the new row 13 so styles (s="4") are reused.
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "给定一个整数数组 nums 和一个目标值 target，请你在该数组中找出和为目标值的那 两个 整数，并返回他们的数组下标。 `n 你可以假设每种输入只会对应一个答案。但是，数组中同一个元素不能使用两遍。 "
$ws.Range("D13").Value = "1 map的key存储目标值target-当前数字的差值，value存储当前位置`n2 遍历数组，当前数字是否在map中，`n       如果不存在就将target-val作为key，index作为value存入map`n        如果存在，就说明找到了这个数字。`n        返回这两个数字的下标"
$ws.Range("E13").Value = "哈希表`n两数之和"
$ws.Range("F13").Value = "O(n)"
$ws.Range("G13").Value = "O(1)"
$ws.Rows.Item(13).RowHeight = 176

# Scroll / select as the author left it, and make the hash sheet the active tab.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("D13").Select()
$ws.Activate()
